$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.859.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "'2.602.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'307.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("D6").Value = "'98.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("D7").Value = "'0.593"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "'38.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'53.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").Value = "'8.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").Value = "'2.996.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "'2.611.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.12%  "
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'14.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").Value = "'45.898.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'12.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.39%  "
$ws.Range("D23").Value = "'285.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.82%  "
$ws.Range("D24").Value = "'73.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'29.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("D29").Value = "'4.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "'10.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'38.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("E33").Value = "  +2.06%  "
$ws.Range("D34").Value = "'3.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").Value = "'158.49"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("D36").Value = "'2.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("E37").Value = "  -3.05%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  +3.35%  "
$ws.Range("D40").Value = "'0.123"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").Value = "'15.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.12%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").Value = "'3.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'4.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").Value = "'21.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'2.113.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("D47").Value = "'0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").Value = "'94.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.32%  "
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").Value = "'108.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("D51").Value = "'2.847.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
